$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header "weight" in column J (row 1), matching existing header style.
$ws.Cells.Item(1, 10).Value = "weight"
$ws.Cells.Item(1, 10).Style = $ws.Cells.Item(1, 9).Style

# Fill weight = 1 for each data row (2 through 16).
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = 1
}

# Move the active selection, as captured in the saved view state.
$ws.Range("K5").Select()
